$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "948×9=8532"; New = "479×2=958" },
    @{ Old = "246×7=1722"; New = "715×4=2860" },
    @{ Old = "848×3=2544"; New = "481×3=1443" },
    @{ Old = "755×6=4530"; New = "671×2=1342" },
    @{ Old = "490×4=1960"; New = "972×7=6804" },
    @{ Old = "970×7=6790"; New = "285×7=1995" },
    @{ Old = "898×6=5388"; New = "403×6=2418" },
    @{ Old = "686×6=4116"; New = "736×5=3680" },
    @{ Old = "748×6=4488"; New = "412×2=824" },
    @{ Old = "817×3=2451"; New = "126×4=504" },
    @{ Old = "620×7=4340"; New = "139×3=417" },
    @{ Old = "249×3=747"; New = "528×6=3168" },
    @{ Old = "480×3=1440"; New = "806×9=7254" },
    @{ Old = "577×2=1154"; New = "102×8=816" },
    @{ Old = "769×3=2307"; New = "721×4=2884" },
    @{ Old = "408×3=1224"; New = "385×4=1540" },
    @{ Old = "864×3=2592"; New = "519×9=4671" },
    @{ Old = "892×5=4460"; New = "927×8=7416" },
    @{ Old = "135×7=945"; New = "320×3=960" },
    @{ Old = "203×3=609"; New = "249×2=498" },
    @{ Old = "990×7=6930"; New = "631×9=5679" },
    @{ Old = "121×7=847"; New = "517×4=2068" },
    @{ Old = "444×8=3552"; New = "550×9=4950" },
    @{ Old = "610×4=2440"; New = "353×4=1412" },
    @{ Old = "668×5=3340"; New = "770×4=3080" }
)

foreach ($pair in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $null = $find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}

Write-Output "Replacements applied: $($replacements.Count)"
